$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.314.65'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.38%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.567.24'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -3.53%  '
$ws.Range('E4').Value = '  -0.40%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '207.44'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.95%  '
$ws.Range('E6').Value = '  -0.39%  '
$ws.Range('E7').Value = '  -4.79%  '
$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.243'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.75%  '
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0607'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.42%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '17.83'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.08%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0782'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.73%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.784.56'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.55%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.568.13'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.33%  '
$ws.Range('E14').Value = '  -3.90%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.507'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.19%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '25.309.64'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.38%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '59.53'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.71%  '
$ws.Range('E18').Value = '  -2.93%  '
$ws.Range('E19').Value = '  -0.27%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '185.59'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.16%  '
$ws.Range('E21').Value = '  -2.17%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.30'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.81%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.89'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.92%  '
$ws.Range('E24').Value = '  -0.41%  '
$ws.Range('E25').Value = '  -3.92%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '140.49'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.08%  '
$ws.Range('E27').Value = '  -6.94%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.45'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.68%  '
$ws.Range('E29').Value = '  -2.01%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0465'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.59%  '
$ws.Range('E32').Value = '  -2.58%  '
$ws.Range('E33').Value = '  -3.46%  '
$ws.Range('E34').Value = '  -1.85%  '
$ws.Range('E35').Value = '  -3.47%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.091.47'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.42%  '
$ws.Range('E38').Value = '  -5.27%  '
$ws.Range('E39').Value = '  -2.42%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.820'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +6.87%  '
$ws.Range('B41').Value = 'ImmutableX'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.495'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.47%  '
$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.771'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -8.83%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '93.28'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.69%  '
$ws.Range('E44').Value = '  -1.95%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.698.95'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.52%  '
$ws.Range('E46').Value = '  -2.73%  '
$ws.Range('E47').Value = '  -3.28%  '
$ws.Range('E48').Value = '  -4.66%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.44'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.68%  '
$ws.Range('E50').Value = '  -1.71%  '
$ws.Range('E51').Value = '  -0.52%  '
